# Remove the "ECs" sending-cluster block (rows 2-4) now that ECs is no longer
# a sending cluster in the updated TPM run, leaving only the FAPs/MuSCs
# sending-cluster rows (which shift up to rows 2-7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:4").Delete()

# Refresh the remaining rows' numeric columns with the values recomputed
# from the new TPM data (ligand/receptor specificity and edge-weight
# columns change because the "ECs" sender was dropped from the specificity
# denominator; receptor-side average/total/specificity for the ECs target
# cluster also reflect the new TPM values).
$ws.Range("I2").Value = 0.8906505749177925
$ws.Range("J2").Value = 0.8906505749177924
$ws.Range("M2").Value = 61.58332300000001
$ws.Range("N2").Value = 184.749969
$ws.Range("O2").Value = 0.9887285514905405
$ws.Range("P2").Value = 0.9887285514905406
$ws.Range("Q2").Value = 797.7870903302824
$ws.Range("R2").Value = 7180.083812972542
$ws.Range("S2").Value = 0.8806116528226862
$ws.Range("T2").Value = 0.8806116528226862
$ws.Range("I3").Value = 0.8906505749177925
$ws.Range("J3").Value = 0.8906505749177924
$ws.Range("O3").Value = 0.004138253805436213
$ws.Range("P3").Value = 0.004138253805436213
$ws.Range("S3").Value = 0.003685738130967506
$ws.Range("T3").Value = 0.003685738130967506
$ws.Range("I4").Value = 0.8906505749177925
$ws.Range("J4").Value = 0.8906505749177924
$ws.Range("O4").Value = 0.007133194704023267
$ws.Range("P4").Value = 0.007133194704023268
$ws.Range("S4").Value = 0.006353183964138876
$ws.Range("T4").Value = 0.006353183964138876
$ws.Range("I5").Value = 0.1093494250822076
$ws.Range("J5").Value = 0.1093494250822076
$ws.Range("M5").Value = 61.58332300000001
$ws.Range("N5").Value = 184.749969
$ws.Range("O5").Value = 0.9887285514905405
$ws.Range("P5").Value = 0.9887285514905406
$ws.Range("Q5").Value = 97.94813153707969
$ws.Range("R5").Value = 881.5331838337172
$ws.Range("S5").Value = 0.1081168986678545
$ws.Range("T5").Value = 0.1081168986678545
$ws.Range("I6").Value = 0.1093494250822076
$ws.Range("J6").Value = 0.1093494250822076
$ws.Range("O6").Value = 0.004138253805436213
$ws.Range("P6").Value = 0.004138253805436213
$ws.Range("S6").Value = 0.0004525156744687076
$ws.Range("T6").Value = 0.0004525156744687076
$ws.Range("I7").Value = 0.1093494250822076
$ws.Range("J7").Value = 0.1093494250822076
$ws.Range("O7").Value = 0.007133194704023267
$ws.Range("P7").Value = 0.007133194704023268
$ws.Range("S7").Value = 0.0007800107398843921
$ws.Range("T7").Value = 0.0007800107398843922
